$wb = $excel.ActiveWorkbook

$newLabel = "31/01/2022 - 06/02/2022"

# Add the new "31/01/2022 - 06/02/2022" weekly row (row 14) to every sheet.

$wsInfanzia = $wb.Worksheets.Item("Infanzia")
$wsInfanzia.Cells.Item(14, 1).Value = $newLabel
$wsInfanzia.Cells.Item(14, 2).Value = 33
$wsInfanzia.Cells.Item(14, 3).Value = 143
$wsInfanzia.Cells.Item(14, 4).Value = 176

$wsPrimaria = $wb.Worksheets.Item("Primaria")
$wsPrimaria.Cells.Item(14, 1).Value = $newLabel
$wsPrimaria.Cells.Item(14, 2).Value = 27
$wsPrimaria.Cells.Item(14, 3).Value = 388
$wsPrimaria.Cells.Item(14, 4).Value = 415

$wsMedia = $wb.Worksheets.Item("Media")
$wsMedia.Cells.Item(14, 1).Value = $newLabel
$wsMedia.Cells.Item(14, 2).Value = 4
$wsMedia.Cells.Item(14, 3).Value = 228
$wsMedia.Cells.Item(14, 4).Value = 232

$wsSuperiore = $wb.Worksheets.Item("Superiore")
$wsSuperiore.Cells.Item(14, 1).Value = $newLabel
$wsSuperiore.Cells.Item(14, 2).Value = 13
$wsSuperiore.Cells.Item(14, 3).Value = 267
$wsSuperiore.Cells.Item(14, 4).Value = 280

$wsTotale = $wb.Worksheets.Item("Totale casi")
$wsTotale.Cells.Item(14, 1).Value = $newLabel
$wsTotale.Cells.Item(14, 2).Value = 77
$wsTotale.Cells.Item(14, 3).Value = 1026
$wsTotale.Cells.Item(14, 4).Value = 1103

# Restore / update the selected cell on each sheet, leaving "Totale casi"
# (the last-saved active tab) selected last so it remains the active sheet.
$wsInfanzia.Range("D15").Select()
$wsPrimaria.Range("E14").Select()
$wsMedia.Range("E14").Select()
$wsSuperiore.Range("E14").Select()
$wsTotale.Range("B15").Select()
